# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These D-column cells get purely-numeric-looking text (e.g. '1.00', '160.54').
# Force Text format first so Excel stores them as strings (matching the source
# data, which keeps these as text) instead of auto-converting to numbers.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D18", "D20", "D21", "D22", "D24", "D25", "D26", "D29", "D30", "D31", "D32", "D33", "D34", "D37", "D38", "D40", "D41", "D42", "D44", "D45", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Cell value updates (Coin / Link / Price / Volume(1h)).
$ws.Range("D2").Value = '66.133.75'
$ws.Range("E2").Value = '  +7.69%  '
$ws.Range("D3").Value = '3.025.11'
$ws.Range("E3").Value = '  +5.48%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '585.88'
$ws.Range("E5").Value = '  +4.19%  '
$ws.Range("D6").Value = '160.54'
$ws.Range("E6").Value = '  +13.36%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.021.86'
$ws.Range("E8").Value = '  +5.45%  '
$ws.Range("D9").Value = '0.518'
$ws.Range("E9").Value = '  +3.87%  '
$ws.Range("D10").Value = '6.98'
$ws.Range("E10").Value = '  +1.67%  '
$ws.Range("D11").Value = '0.155'
$ws.Range("E11").Value = '  +6.06%  '
$ws.Range("D12").Value = '0.454'
$ws.Range("E12").Value = '  +6.22%  '
$ws.Range("E13").Value = '  +9.80%  '
$ws.Range("D14").Value = '34.62'
$ws.Range("E14").Value = '  +10.01%  '
$ws.Range("E15").Value = '  +0.63%  '
$ws.Range("D16").Value = '66.140.58'
$ws.Range("E16").Value = '  +7.67%  '
$ws.Range("D17").Value = '3.526.92'
$ws.Range("E17").Value = '  +5.46%  '
$ws.Range("D18").Value = '6.97'
$ws.Range("E18").Value = '  +7.44%  '
$ws.Range("D19").Value = '3.022.05'
$ws.Range("E19").Value = '  +5.15%  '
$ws.Range("D20").Value = '459.66'
$ws.Range("E20").Value = '  +7.62%  '
$ws.Range("D21").Value = '13.86'
$ws.Range("E21").Value = '  +7.06%  '
$ws.Range("D22").Value = '0.685'
$ws.Range("E22").Value = '  +6.07%  '
$ws.Range("E23").Value = '  +9.03%  '
$ws.Range("D24").Value = '82.44'
$ws.Range("E24").Value = '  +4.89%  '
$ws.Range("D25").Value = '2.26'
$ws.Range("E25").Value = '  +12.92%  '
$ws.Range("D26").Value = '12.45'
$ws.Range("E26").Value = '  +5.18%  '
$ws.Range("E27").Value = '  +5.74%  '
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").Value = '8.08'
$ws.Range("E29").Value = '  +16.10%  '
$ws.Range("D30").Value = '2.37'
$ws.Range("E30").Value = '  +18.77%  '
$ws.Range("D31").Value = '0.0000103'
$ws.Range("E31").Value = '  -3.85%  '
$ws.Range("D32").Value = '2.60'
$ws.Range("E32").Value = '  +5.18%  '
$ws.Range("D33").Value = '26.96'
$ws.Range("E33").Value = '  +6.01%  '
$ws.Range("D34").Value = '0.110'
$ws.Range("E34").Value = '  +5.11%  '
$ws.Range("E35").Value = '  -0.33%  '
$ws.Range("E36").Value = '  +5.14%  '
$ws.Range("D37").Value = '5.76'
$ws.Range("E37").Value = '  +8.57%  '
$ws.Range("D38").Value = '2.16'
$ws.Range("E38").Value = '  +14.88%  '
$ws.Range("E39").Value = '  +10.64%  '
$ws.Range("D40").Value = '49.84'
$ws.Range("E40").Value = '  +2.18%  '
$ws.Range("B41").Value = 'Arweave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D41").Value = '44.54'
$ws.Range("E41").Value = '  +14.87%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '0.121'
$ws.Range("E42").Value = '  +7.41%  '
$ws.Range("E43").Value = '  +13.54%  '
$ws.Range("D44").Value = '8.44'
$ws.Range("E44").Value = '  +4.07%  '
$ws.Range("D45").Value = '385.21'
$ws.Range("E45").Value = '  +12.65%  '
$ws.Range("D46").Value = '2.804.01'
$ws.Range("E46").Value = '  +4.90%  '
$ws.Range("E47").Value = '  +6.46%  '
$ws.Range("D48").Value = '134.69'
$ws.Range("E48").Value = '  +2.82%  '
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").Value = '23.63'
$ws.Range("E50").Value = '  +11.68%  '
$ws.Range("D51").Value = '0.106'
$ws.Range("E51").Value = '  +4.38%  '
